# Apply "feat: add 2022-Q1 data" change:
#  1. Insert a new worksheet "2022-Q1" (fund holding detail) right before "总计".
#  2. Prepend a "2022-Q1" summary row to the "总计" sheet, renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: style a cell like the workbook's "header/index" look
#   (bold font, thin box border, centered horizontally, top vertically)
# ---------------------------------------------------------------------------
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# Helper: write a value as genuine TEXT (leading apostrophe forces text even
# when the value looks numeric, e.g. "19.70" / "010714").
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
}

# ---------------------------------------------------------------------------
# 1) New sheet "2022-Q1" inserted right before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row
Set-TextValue $q1.Cells.Item(1,2) "基金代码"
Set-TextValue $q1.Cells.Item(1,3) "基金名称"
Set-TextValue $q1.Cells.Item(1,4) "基金规模"
Set-TextValue $q1.Cells.Item(1,5) "股票总仓位"
Set-TextValue $q1.Cells.Item(1,6) "仓位占比"
Set-TextValue $q1.Cells.Item(1,7) "持有市值(亿元)"
Set-TextValue $q1.Cells.Item(1,8) "仓位排名"
for ($c = 2; $c -le 8; $c++) {
    Set-HeaderStyle $q1.Cells.Item(1, $c)
}

# Data rows: index, code, name, scale, stock-position, ratio, held-value, rank
$rows = @(
    @(0, "010714", "东方红远见价值混合",           "19.70", "86.34", "4.06", "0.7998", 4),
    @(1, "012366", "上投摩根安荣回报混合型证券投资基金A", "41.47", "21.90", "0.90", "0.3732", 10),
    @(2, "004738", "上投摩根安隆回报混合A",          "23.04", "21.18", "1.03", "0.2373", 9),
    @(3, "012367", "上投摩根安荣回报混合型证券投资基金C", "23.94", "21.90", "0.90", "0.2155", 10),
    @(4, "004739", "上投摩根安隆回报混合C",          "7.32",  "21.18", "1.03", "0.0754", 9),
    @(5, "710301", "富安达增强收益债券A",            "0.61",  "20.20", "1.56", "0.0095", 5),
    @(6, "710302", "富安达增强收益债券C",            "0.26",  "20.20", "1.56", "0.0041", 5)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    Set-HeaderStyle $idxCell

    Set-TextValue $q1.Cells.Item($r, 2) $row[1]
    Set-TextValue $q1.Cells.Item($r, 3) $row[2]
    Set-TextValue $q1.Cells.Item($r, 4) $row[3]
    Set-TextValue $q1.Cells.Item($r, 5) $row[4]
    Set-TextValue $q1.Cells.Item($r, 6) $row[5]
    Set-TextValue $q1.Cells.Item($r, 7) $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" sheet and renumber the index column
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

$a2 = $ws.Cells.Item(2, 1)
$a2.Value = 0
Set-HeaderStyle $a2

$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 7
$ws.Cells.Item(2, 4).Value = 1.71
$ws.Range("B2:D2").ClearFormats()

# Renumber the index column (A) for the rows that got pushed down
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4
